$d = $word.ActiveDocument

$pairs = @(
    @{old="71×30="; new="62×72="},
    @{old="13×64="; new="24×30="},
    @{old="36×43="; new="12×15="},
    @{old="95×25="; new="21×56="},
    @{old="35×41="; new="48×43="},
    @{old="73×71="; new="30×51="},
    @{old="96×14="; new="33×42="},
    @{old="39×41="; new="65×89="},
    @{old="86×30="; new="15×65="},
    @{old="33×34="; new="80×54="},
    @{old="14×64="; new="46×24="},
    @{old="89×35="; new="30×41="},
    @{old="86×51="; new="63×68="},
    @{old="95×92="; new="81×25="},
    @{old="20×13="; new="27×14="},
    @{old="32×17="; new="52×24="},
    @{old="80×70="; new="63×38="},
    @{old="34×52="; new="82×89="},
    @{old="28×83="; new="63×39="},
    @{old="14×61="; new="42×46="},
    @{old="91×45="; new="68×67="},
    @{old="28×41="; new="77×99="},
    @{old="97×21="; new="96×93="},
    @{old="74×50="; new="52×89="},
    @{old="39×79="; new="80×12="}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}

$d.Save()
